# CustCross.xlsx — update the two multi-line remark/note cells on the
# "DBD" sheet (row 10 = SubCompanyCode remark, row 11 = CrossUse remark):
# the "XX: " colon+space separators in the enumerated notes become "XX:"
# (space after the colon removed), and leave the cursor/active sheet on
# "DBD" at G10 (mirrors the author re-selecting DBD after editing G10/G11,
# instead of DBS being the active tab as before).

$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")

$wsDBD.Range("G10").Value = "共用代碼檔`n01:新光金控`n02:新光人壽`n03:新光銀行`n04:新光信託`n05:保險經紀人`n06:元富證券"
$wsDBD.Range("G11").Value = "Y:同意使用`nN:不同意使用"

$wsDBD.Activate() | Out-Null
$wsDBD.Range("G10").Select() | Out-Null
